$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$listBulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

function New-Para([string]$innerXml) {
    return "<w:p $wNs>$innerXml</w:p>"
}

# ---------------------------------------------------------------------
# 1. Shorten the title and remove the separate "Meta description" line.
#    Both paragraphs get collapsed into a single Heading1 paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$metaPara  = $d.Paragraphs.Item(2)
$combined  = $d.Range($titlePara.Range.Start, $metaPara.Range.End)
$xml = New-Para '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Da Vinci Diamonds Dual Play Free</w:t></w:r>'
$combined.InsertXML($xml)

# ---------------------------------------------------------------------
# 2. Update the "What we like" bullet list (4 consecutive paragraphs).
#    (Step 1 merged two paragraphs into one, so every later paragraph
#    index is now one less than in the original document.)
# ---------------------------------------------------------------------
$firstLike = $d.Paragraphs.Item(47)
$lastLike  = $d.Paragraphs.Item(50)
$combined  = $d.Range($firstLike.Range.Start, $lastLike.Range.End)

$xml  = New-Para "$listBulletPPr<w:r/><w:r><w:t>Tumbling Reels mode for increased chances of winning</w:t></w:r>"
$xml += New-Para "$listBulletPPr<w:r/><w:r><w:t>Wild symbol and Extra Paylines Bonus feature</w:t></w:r>"
$xml += New-Para "$listBulletPPr<w:r/><w:r><w:t>Dual game panels for double the chances of winning</w:t></w:r>"
$xml += New-Para "$listBulletPPr<w:r/><w:r><w:t>Luxurious design featuring Da Vinci artwork</w:t></w:r>"
$combined.InsertXML($xml)

# ---------------------------------------------------------------------
# 3. Update the "What we don't like" bullet list (2 consecutive paragraphs).
# ---------------------------------------------------------------------
$firstDontLike = $d.Paragraphs.Item(52)
$lastDontLike  = $d.Paragraphs.Item(53)
$combined      = $d.Range($firstDontLike.Range.Start, $lastDontLike.Range.End)

$xml  = New-Para "$listBulletPPr<w:r/><w:r><w:t>No progressive jackpot feature</w:t></w:r>"
$xml += New-Para "$listBulletPPr<w:r/><w:r><w:t>Limited bonus features</w:t></w:r>"
$combined.InsertXML($xml)

# ---------------------------------------------------------------------
# 4. Insert a new bold "Play Da Vinci Diamonds Dual Play Free" paragraph
#    and rewrite the closing italic paragraph's text.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$xml  = New-Para '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Da Vinci Diamonds Dual Play Free</w:t></w:r>'
$xml += New-Para '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Da Vinci Diamonds Dual Play slot game and play for free.</w:t></w:r>'
$lastPara.Range.InsertXML($xml)
